# Update the lognormal distribution parameter strings (format change: add
# a second parameter, e.g. a standard deviation, to each "lognorm,<mean>"
# string, and correct a stray "lognorm,6" typo to "lognorm,0.6,0.06") on
# every trajectory worksheet (traj1..traj5), then re-select the cell the
# change was made in on each sheet, matching the author's final edit.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("traj1")
$ws1.Range("B2").Value = "lognorm,0.6,0.06"
$ws1.Range("B2").Select()

$ws2 = $wb.Worksheets.Item("traj2")
$ws2.Range("B2").Value = "lognorm,0.6,0.06"
$ws2.Range("B2").Select()

$ws3 = $wb.Worksheets.Item("traj3")
$ws3.Range("B2").Value = "lognorm,0.6,0.06"
$ws3.Range("B3").Value = "lognorm,5.4,0.54"
$ws3.Range("B2").Select()

$ws4 = $wb.Worksheets.Item("traj4")
$ws4.Range("B2").Value = "lognorm,5.4,0.54"
$ws4.Range("B3").Value = "lognorm,0.6,0.06"
$ws4.Range("B3").Select()

$ws5 = $wb.Worksheets.Item("traj5")
$ws5.Range("B2").Value = "lognorm,0.6,0.06"
$ws5.Range("B3").Value = "lognorm,4.8,0.48"
$ws5.Range("B4").Value = "lognorm,0.6,0.06"
$ws5.Range("B4").Select()

$ws1.Select()
